# Generate Report for Handoff
#
# Refresh the localization-status report with the newly generated source
# file (UUID rotated from 0bf25ae4-... to 2568fbfd-...), its new handoff
# xliff hashes/filenames, and the refreshed "ready for handoff" timestamps.

$wb = $excel.ActiveWorkbook

$oldId = "0bf25ae4-b59e-46d0-9d97-fcce1a4d0632"
$newId = "2568fbfd-1eb4-41a6-8c96-17490e131c4f"

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Overview sheet ------------------------------------------------------
$ws1.Range("A2").Value = "$newId.md"

# B2 carries a hyperlink - update its display text (and therefore the cell
# value) through the hyperlink object itself so the existing hyperlink
# relationship/formatting is left untouched.
$ws1.Range("B2").Hyperlinks.Item(1).TextToDisplay = "e2e\$newId.md"

$ws1.Range("G2").Value = "2016-08-29 17:04:35"

# --- zh-cn sheet -----------------------------------------------------------
$ws2.Range("A2").Hyperlinks.Item(1).TextToDisplay = "$newId.md"
$ws2.Range("G2").Value = "$newId.c86aa87242498da3292227c9e50b18db13db438f.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-29 17:04:30"

# --- de-de sheet -----------------------------------------------------------
$ws3.Range("A2").Hyperlinks.Item(1).TextToDisplay = "$newId.md"
$ws3.Range("G2").Value = "$newId.c86aa87242498da3292227c9e50b18db13db438f.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-29 17:04:35"
